$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, centered, bordered) from H1 onto the
# two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2..49 (columns I = "I0", J = "IF")
$iVals = @(5,9,5,7,6,9,4,5,5,10,6,6,8,5,8,5,8,6,2,9,7,6,6,6,9,5,5,5,2,7,9,6,10,5,7,9,6,7,5,5,9,8,7,4,7,7,7,6)
$jVals = @(5,9,5,8,6,9,5,6,6,10,7,7,9,6,8,6,9,7,4,9,8,7,7,7,9,5,6,5,5,8,9,7,10,5,7,9,6,7,5,5,9,8,7,5,7,7,7,6)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
